$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'modest medical',
    'moisture wicking compression leggings',
    'motion control knee brace',
    'moto leggings',
    'moving runners',
    'mud light football',
    'mud pants',
    'mud riding games',
    'mud run gear',
    'mud runner game',
    'muffin time',
    'muffin top pants',
    'muffin top workout',
    'mujer fantastica',
    'multicolor puff ball',
    'muscle amp x',
    'muscle and joint',
    'muscle and joint pain',
    'muscle body wash',
    'muscle building equipment',
    'muscle building for women',
    'muscle chart',
    'muscle club apparel',
    'muscle compression',
    'muscle compression bands',
    'muscle endurance',
    'muscle exercise',
    'muscle fatigue',
    'muscle flexor',
    'muscle for men',
    'muscle gear',
    'muscle gear clothing',
    'muscle joint',
    'muscle leggings',
    'muscle leggings for women',
    'muscle leggings men',
    'muscle mil',
    'muscle pain tape',
    'muscle pants',
    'muscle pro',
    'muscle recovery equipment',
    'muscle short',
    'muscle stretch',
    'muscle stretching',
    'muscle support',
    'muscle support tape',
    'muscle tear',
    'muscle tech',
    'muscle tech multi',
    'muscle tech muscle',
    'muscle tech power',
    'muscle tightening',
    'muscle top women',
    'muscle training',
    'muscle training gear',
    'muscle training gear abs',
    'muscle up rings',
    'muscle wear men',
    'muscle women',
    'muscle x',
    'muscular model',
    'must haves women',
    'my fitness',
    'n tight',
    'naked women',
    'nasty sounds',
    'nasty woman club',
    'natori tights',
    'natural grip',
    'natural joint support',
    'natural performance for women',
    'natural training pants',
    'natures way bone',
    'navy blue capri tights',
    'navy blue compression leggings women',
    'navy blue compression tights for women',
    'navy blue leggins womens',
    'navy blue maternity leggings',
    'navy capris',
    'navy compression leggings',
    'navy leggins',
    'neoprene capri',
    'neoprene capri pants',
    'neoprene capri women',
    'neoprene clothing for women',
    'neoprene compression',
    'neoprene compression shirt',
    'neoprene crop',
    'neoprene hose',
    'neoprene knee supports',
    'neoprene leggings',
    'neoprene leggings women',
    'neoprene men',
    'neoprene pants',
    'neoprene seam tape',
    'neoprene sweat pants',
    'neoprene thigh support',
    'neoprene weightlifting',
    'nerve compression',
    'nerve gear'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

